$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "body" column (E), pushing it to F.
$ws.Columns("E").Insert()

# New header + value for the inserted "senderEmail" column.
$ws.Range("E1").Value = "senderEmail"
$ws.Range("E2").Value = "zephytk@gmail.com"

# Turn E2 into a mailto hyperlink, then restore the clean Hyperlink cell style
# (copying the format from the existing hyperlink cell B2 avoids leaving the
# cell with a combined/duplicate style).
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:zephytk@gmail.com")
$ws.Range("B2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

# Update the body text (now in column F) to add the extra "Workflow Version 5" line.
$ws.Range("F2").Value = "This is the start of email body message`nWorkflow Version 5`nTeam Charlie Rocks`nThis is  a bot generated message.`nregards"

# Undo the automatic row-height expansion caused by the taller wrapped text.
$ws.Rows(2).AutoFit()

# Match column E's width to column D's width.
$ws.Columns("E").ColumnWidth = $ws.Columns("D").ColumnWidth()

# Update the active selection.
$ws.Range("E8").Select()

Write-Host "done"
